$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.071684527833219
$ws.Range("D2").Value = 1.073303175105101
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.083105220937883
$ws.Range("I2").Value = 1.061812187877856
$ws.Range("J2").Value = 1.076607878401696
$ws.Range("K2").Value = 1.075995393738399
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.085771702332616
$ws.Range("N2").Value = 1.078136785077677

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.072914509975875
$ws.Range("D3").Value = 1.074287236733825
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.084229850560326
$ws.Range("I3").Value = 1.062262765656662
$ws.Range("J3").Value = 1.077494556110752
$ws.Range("K3").Value = 1.076795953142709
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.086714349908488
$ws.Range("N3").Value = 1.079024721970782

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.073709616925371
$ws.Range("D4").Value = 1.074923197423882
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.084957094286641
$ws.Range("I4").Value = 1.062552616829222
$ws.Range("J4").Value = 1.078066962681163
$ws.Range("K4").Value = 1.077312557013944
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.087323245960328
$ws.Range("N4").Value = 1.079597941424179

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.074043697395887
$ws.Range("D5").Value = 1.075190367086041
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.085262717629916
$ws.Range("I5").Value = 1.062674063692369
$ws.Range("J5").Value = 1.078307285094239
$ws.Range("K5").Value = 1.077529400981081
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.087578973846526
$ws.Range("N5").Value = 1.079838605122646

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.074099780355229
$ws.Range("D6").Value = 1.075235215058377
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.085314026741824
$ws.Range("I6").Value = 1.062694431344948
$ws.Range("J6").Value = 1.078347617731708
$ws.Range("K6").Value = 1.077565790364056
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.087621896956231
$ws.Range("N6").Value = 1.079878995037088

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.073714081638604
$ws.Range("D7").Value = 1.074926768095495
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.084961178470318
$ws.Range("I7").Value = 1.062554241204349
$ws.Range("J7").Value = 1.078070175125496
$ws.Range("K7").Value = 1.077315455813034
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.087326663996929
$ws.Range("N7").Value = 1.079601158430552

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.07210036791965
$ws.Range("D8").Value = 1.07363590851626
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.083485392558458
$ws.Range("I8").Value = 1.061964815970072
$ws.Range("J8").Value = 1.076907812454246
$ws.Range("K8").Value = 1.076266239805085
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.0860904950177
$ws.Range("N8").Value = 1.078437145070984

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.069250748940058
$ws.Range("D9").Value = 1.07135510213928
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.080881210689818
$ws.Range("I9").Value = 1.060913083539198
$ws.Range("J9").Value = 1.074849286449386
$ws.Range("K9").Value = 1.074406502245883
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.083904010752167
$ws.Range("N9").Value = 1.07637569572308

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.067346761157851
$ws.Range("D10").Value = 1.06983032643094
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.079142510844312
$ws.Range("I10").Value = 1.060203058338178
$ws.Range("J10").Value = 1.073469895099865
$ws.Range("K10").Value = 1.073159253641466
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.08244073152349
$ws.Range("N10").Value = 1.074994345479625

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.066521267759813
$ws.Range("D11").Value = 1.069169052275975
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.078388995749218
$ws.Range("I11").Value = 1.059893489792527
$ws.Range("J11").Value = 1.072870906247011
$ws.Range("K11").Value = 1.072617396406359
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.081805756066869
$ws.Range("N11").Value = 1.074394505993895

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.066214480877831
$ws.Range("D12").Value = 1.06892326774566
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.078109007098555
$ws.Range("I12").Value = 1.059778181859622
$ws.Range("J12").Value = 1.072648156774546
$ws.Range("K12").Value = 1.072415855356138
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.081569690387624
$ws.Range("N12").Value = 1.074171440191628

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.06628029508727
$ws.Range("D13").Value = 1.068975996540745
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.078169070160898
$ws.Range("I13").Value = 1.059802930324264
$ws.Range("J13").Value = 1.072695949022388
$ws.Range("K13").Value = 1.072459098900395
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.081620336705228
$ws.Range("N13").Value = 1.074219300309945

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.066495911983063
$ws.Range("D14").Value = 1.069148738883128
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.078365853845918
$ws.Range("I14").Value = 1.059883964947074
$ws.Range("J14").Value = 1.072852498990395
$ws.Range("K14").Value = 1.072600742504147
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.081786247050982
$ws.Range("N14").Value = 1.074376072596863

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.066628739113122
$ws.Range("D15").Value = 1.069255150157124
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.078487085498231
$ws.Range("I15").Value = 1.059933850554443
$ws.Range("J15").Value = 1.072948920286388
$ws.Range("K15").Value = 1.072687977808722
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.081888442322671
$ws.Range("N15").Value = 1.074472630822156

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06740152392762
$ws.Range("D16").Value = 1.069874191000551
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.079192505327676
$ws.Range("I16").Value = 1.060223558531928
$ws.Range("J16").Value = 1.073509611922648
$ws.Range("K16").Value = 1.073195177016135
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.082482843778354
$ws.Range("N16").Value = 1.075034118704852

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.067885986732622
$ws.Range("D17").Value = 1.070262220002514
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.079634821712041
$ws.Range("I17").Value = 1.060404715301682
$ws.Range("J17").Value = 1.07386086123108
$ws.Range("K17").Value = 1.073512848822317
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.082855328789231
$ws.Range("N17").Value = 1.075385866827591

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.068168464221405
$ws.Range("D18").Value = 1.070488450924395
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.079892755156585
$ws.Range("I18").Value = 1.060510176245139
$ws.Range("J18").Value = 1.074065574833041
$ws.Range("K18").Value = 1.073697968816236
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.083072461362834
$ws.Range("N18").Value = 1.075590871146348

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.068264764639822
$ws.Range("D19").Value = 1.070565572955099
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.079980693315678
$ws.Range("I19").Value = 1.060546101018147
$ws.Range("J19").Value = 1.074135349107107
$ws.Range("K19").Value = 1.073761060739983
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.083146475701666
$ws.Range("N19").Value = 1.075660744507885

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.067834018980651
$ws.Range("D20").Value = 1.070220598505596
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.0795873718045
$ws.Range("I20").Value = 1.060385300082537
$ws.Range("J20").Value = 1.073823192504741
$ws.Range("K20").Value = 1.073478783512939
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.082815378305712
$ws.Range("N20").Value = 1.075348144607339

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.066432422678256
$ws.Range("D21").Value = 1.069097874944275
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.078307908724415
$ws.Range("I21").Value = 1.059860111126377
$ws.Range("J21").Value = 1.0728064060779
$ws.Range("K21").Value = 1.072559039495635
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.081737396336618
$ws.Range("N21").Value = 1.074329914227146

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.065550245125862
$ws.Range("D22").Value = 1.068391059626448
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.077502882061367
$ws.Range("I22").Value = 1.059528049580394
$ws.Range("J22").Value = 1.072165616188799
$ws.Range("K22").Value = 1.071979190076476
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.081058424475051
$ws.Range("N22").Value = 1.073688214342902

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.066017994255807
$ws.Range("D23").Value = 1.068765843212785
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.077929697483894
$ws.Range("I23").Value = 1.059704257948119
$ws.Range("J23").Value = 1.072505453646007
$ws.Range("K23").Value = 1.072286728754109
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.081418474998401
$ws.Range("N23").Value = 1.07402853440828

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.067857501284274
$ws.Range("D24").Value = 1.070239405776944
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.079608812567514
$ws.Range("I24").Value = 1.060394073616925
$ws.Range("J24").Value = 1.073840213887983
$ws.Range("K24").Value = 1.073494176693191
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.082833430616271
$ws.Range("N24").Value = 1.075365190162898

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.069988178557669
$ws.Range("D25").Value = 1.071945484503146
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.08155489987581
$ws.Range("I25").Value = 1.061186539922608
$ws.Range("J25").Value = 1.075382696958984
$ws.Range("K25").Value = 1.074888589700419
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.084470252194629
$ws.Range("N25").Value = 1.076909863736785
